$d = $word.ActiveDocument

$replacements = @(
    @{old="59×45=2655"; new="69×90=6210"},
    @{old="51×38=1938"; new="97×44=4268"},
    @{old="85×36=3060"; new="13×40=520"},
    @{old="91×49=4459"; new="63×25=1575"},
    @{old="52×97=5044"; new="82×19=1558"},
    @{old="20×22=440"; new="44×18=792"},
    @{old="72×13=936"; new="38×55=2090"},
    @{old="64×46=2944"; new="23×29=667"},
    @{old="85×24=2040"; new="21×28=588"},
    @{old="61×98=5978"; new="76×19=1444"},
    @{old="20×23=460"; new="51×49=2499"},
    @{old="28×93=2604"; new="91×26=2366"},
    @{old="16×86=1376"; new="60×34=2040"},
    @{old="94×85=7990"; new="84×28=2352"},
    @{old="27×69=1863"; new="51×50=2550"},
    @{old="81×75=6075"; new="33×96=3168"},
    @{old="71×55=3905"; new="89×83=7387"},
    @{old="81×73=5913"; new="88×35=3080"},
    @{old="27×19=513"; new="46×59=2714"},
    @{old="33×93=3069"; new="78×84=6552"},
    @{old="13×54=702"; new="91×84=7644"},
    @{old="17×35=595"; new="22×52=1144"},
    @{old="34×28=952"; new="55×62=3410"},
    @{old="49×39=1911"; new="62×84=5208"},
    @{old="15×45=675"; new="66×39=2574"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
